$wb = $excel.ActiveWorkbook
$waitSheet = $wb.Worksheets.Item("Wait Events per SQL Statement")

# Add the new sheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "All Wait Events"

# Copy the "Wait Event" detail header/sample rows (rows 2 and 4 of the
# "Wait Events per SQL Statement" sheet) into the new sheet, reusing their
# cell formatting (font/fill/border/number format).
$waitSheet.Range("A2:D2").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

$waitSheet.Range("A4:D4").Copy()
$ws.Range("A2:D2").PasteSpecial(-4122)

# Match column widths to the "Wait Events per SQL Statement" sheet
$ws.Columns("A").ColumnWidth = 19.833333333333336
$ws.Columns("B:C").ColumnWidth = 29.833333333333336
$ws.Columns("D").ColumnWidth = 22.666666666666664

# Freeze the header row and set the view/zoom to match
$ws.Activate()
$excel.ActiveWindow.Zoom = 120
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("A2").Select()
